# Update the "phone" column (B) on the active sheet:
#   - header text: "טלפון"      -> "מספר טלפון"
#   - B2 value:    507676706     -> "050-7676706" (text, dash formatted)
#   - B3 value:    586208430     -> "058-6208430" (text, dash formatted)
# Because the new values contain a literal dash and a leading zero, the
# whole column is (re)formatted as Text ("@") first so the values are
# stored/typed as text rather than being re-interpreted as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column B ("phone number") as Text so the dashed values stick as
# strings instead of being coerced back into numbers.
$ws.Columns.Item(2).NumberFormat = "@"

# New header label for column B.
$ws.Range("B1").Value = "מספר טלפון"

# New phone number values, written as formatted text strings.
$ws.Range("B2").Value = "050-7676706"
$ws.Range("B3").Value = "058-6208430"

# Leave the same cell selected/active as in the saved workbook.
$ws.Range("J5").Select()
